# Update the "last updated" timestamp banner in A1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 4 de Agosto de 2020 a las 10:16"

# Update country statistics rows (columns: B=Casos totales, C=Nuevos casos,
# D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes)

# Row 7
$ws.Range("B7").Value = 861423
$ws.Range("C7").Value = 5159
$ws.Range("D7").Value = 661471
$ws.Range("E7").Value = 185601
$ws.Range("G7").Value = 144
$ws.Range("H7").Value = 14351

# Row 30
$ws.Range("E30").Value = 27630
$ws.Range("H30").Value = 1058

# Row 46
$ws.Range("B46").Value = 53346
$ws.Range("C46").Value = 295
$ws.Range("E46").Value = 6140

# Row 56
$ws.Range("B56").Value = 37541
$ws.Range("C56").Value = 412
$ws.Range("D56").Value = 28743
$ws.Range("E56").Value = 7371
$ws.Range("G56").Value = 7
$ws.Range("H56").Value = 1427

# Row 106
$ws.Range("B106").Value = 4553
$ws.Range("C106").Value = 9
$ws.Range("D106").Value = 3415
$ws.Range("E106").Value = 540
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 598

# Row 125
$ws.Range("B125").Value = 2368
$ws.Range("C125").Value = 14
$ws.Range("D125").Value = 1771
$ws.Range("E125").Value = 568

# Row 129
$ws.Range("B129").Value = 2091
$ws.Range("C129").Value = 11
$ws.Range("D129").Value = 1937
$ws.Range("E129").Value = 91

# Row 162
$ws.Range("E162").Value = 270
$ws.Range("G162").Value = 2
$ws.Range("H162").Value = 8
